$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A1").Value = "SearchKey"
$ws.Range("A2").Value = "iPhone"

# Excel's stored column <width> includes the ~0.8333 char padding on top of
# the ColumnWidth value you assign, so back it out to land on exactly 18.
$ws.Columns.Item(1).ColumnWidth = 17.1666667

$ws.Range("A4").Select()
